$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16-22: Samay Diaz Jimenez now occupies the first 7 data rows,
# carrying periods 2105-2111 (previously held by the other 7 workers
# with period 2112 for one row each).
$ws.Range("C16").Value = "1001834084"
$ws.Range("D16").Value = "SAMAY DIAZ JIMENEZ"
$ws.Range("E16").Value = "2105"
$ws.Range("F16").Value = 31495

$ws.Range("C17").Value = "1001834084"
$ws.Range("D17").Value = "SAMAY DIAZ JIMENEZ"
$ws.Range("E17").Value = "2106"
$ws.Range("F17").Value = 36341

$ws.Range("C18").Value = "1001834084"
$ws.Range("D18").Value = "SAMAY DIAZ JIMENEZ"
$ws.Range("E18").Value = "2107"
$ws.Range("F18").Value = 36341

$ws.Range("C19").Value = "1001834084"
$ws.Range("D19").Value = "SAMAY DIAZ JIMENEZ"
$ws.Range("E19").Value = "2108"
$ws.Range("F19").Value = 36341

$ws.Range("C20").Value = "1001834084"
$ws.Range("D20").Value = "SAMAY DIAZ JIMENEZ"
$ws.Range("E20").Value = "2109"
$ws.Range("F20").Value = 36341

$ws.Range("C21").Value = "1001834084"
$ws.Range("D21").Value = "SAMAY DIAZ JIMENEZ"
$ws.Range("E21").Value = "2110"
$ws.Range("F21").Value = 36341

$ws.Range("C22").Value = "1001834084"
$ws.Range("D22").Value = "SAMAY DIAZ JIMENEZ"
$ws.Range("E22").Value = "2111"
$ws.Range("F22").Value = 36341

# Rows 23-29: the 7 workers that used to be in rows 16-22 move down here,
# all now on period 2112, and their Salario Basico (G) updates to 877803.
$ws.Range("C23").Value = "45691991"
$ws.Range("D23").Value = "YUDIS MARGOTH TAPIAS GUTIERREZ"
$ws.Range("E23").Value = "2112"
$ws.Range("F23").Value = 18726
$ws.Range("G23").Value = 877803

$ws.Range("C24").Value = "45521649"
$ws.Range("D24").Value = "ALEXANDRA PATRICIA RUIZ MONTIEL"
$ws.Range("E24").Value = "2112"
$ws.Range("F24").Value = 18726
$ws.Range("G24").Value = 877803

$ws.Range("C25").Value = "45516185"
$ws.Range("D25").Value = "YESENIA JULIO SALAS"
$ws.Range("E25").Value = "2112"
$ws.Range("F25").Value = 18726
$ws.Range("G25").Value = 877803

$ws.Range("C26").Value = "1063720644"
$ws.Range("D26").Value = "SANTIAGO BLANCO AGAMEZ"
$ws.Range("E26").Value = "2112"
$ws.Range("F26").Value = 18726
$ws.Range("G26").Value = 877803

$ws.Range("C27").Value = "45423586"
$ws.Range("D27").Value = "AMADA ISABEL JIMENEZ BENAVIDES"
$ws.Range("E27").Value = "2112"
$ws.Range("F27").Value = 18726
$ws.Range("G27").Value = 877803

$ws.Range("C28").Value = "1151448981"
$ws.Range("D28").Value = "KEYLA MARCELA TORRES HERNANDEZ"
$ws.Range("E28").Value = "2112"
$ws.Range("F28").Value = 18726
$ws.Range("G28").Value = 877803

$ws.Range("C29").Value = "1070818231"
$ws.Range("D29").Value = "LUIS ALBERTO BLANCO AGAMEZ"
$ws.Range("E29").Value = "2112"
$ws.Range("F29").Value = 18726
$ws.Range("G29").Value = 877803

# Row 30: Samay Diaz Jimenez's final row also moves to period 2112.
$ws.Range("C30").Value = "1001834084"
$ws.Range("D30").Value = "SAMAY DIAZ JIMENEZ"
$ws.Range("E30").Value = "2112"
$ws.Range("F30").Value = 19382
